# The document's very first paragraph holds nothing but a single run
# containing a floating "Rectangle 1" text box (the "School Logo and
# School information" banner), anchored at the top of the page. The
# commit removes that run entirely so the first paragraph goes back to
# being empty (just its paragraph mark / pPr), fixing the header
# placement at the start of the document.
#
# The second floating shape ("Rectangle 4" - "Scheme of Work Template")
# lives further down in the body and must be left untouched.

$d = $word.ActiveDocument

$shapeCountBefore = $d.Shapes.Count

# Remove the floating shape anchored in the document's first paragraph
# (the "School Logo and School information" rectangle/text box).
$d.Shapes.Item(1).Delete()

$shapeCountAfter = $d.Shapes.Count

Write-Output "Shapes before: $shapeCountBefore, after: $shapeCountAfter"
